$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    # Force the cell to be stored as text (matches the source inlineStr cells),
    # avoiding Excel's automatic "looks like a number" coercion + float noise,
    # then strip the temporary number-format override so no stray style sticks.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextCell "D2" "63.255.69"
$ws.Range("E2").Value = "  +3.37%  "

# Row 3 - Ethereum
Set-TextCell "D3" "3.035.70"
$ws.Range("E3").Value = "  +1.82%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.08%  "

# Row 5 - BNB
Set-TextCell "D5" "595.56"
$ws.Range("E5").Value = "  +0.03%  "

# Row 6 - Solana
Set-TextCell "D6" "154.26"
$ws.Range("E6").Value = "  +7.75%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.05%  "

# Row 8 - LidoStakedEther
Set-TextCell "D8" "3.032.65"
$ws.Range("E8").Value = "  +1.77%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  +0.44%  "

# Row 10 - Toncoin
Set-TextCell "D10" "6.99"
$ws.Range("E10").Value = "  +16.17%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +3.17%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +2.55%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  +4.05%  "

# Row 14 - Avalanche
Set-TextCell "D14" "35.75"
$ws.Range("E14").Value = "  +4.71%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  -0.76%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-TextCell "D16" "3.538.57"
$ws.Range("E16").Value = "  +1.86%  "

# Row 17 - Polkadot
$ws.Range("E17").Value = "  +3.70%  "

# Row 18 - WrappedBTC
Set-TextCell "D18" "63.208.74"
$ws.Range("E18").Value = "  +3.39%  "

# Row 19 - WrappedEther
Set-TextCell "D19" "3.037.13"
$ws.Range("E19").Value = "  +2.06%  "

# Row 20 - BitcoinCash
Set-TextCell "D20" "453.02"
$ws.Range("E20").Value = "  +0.80%  "

# Row 21 - Chainlink
Set-TextCell "D21" "14.29"
$ws.Range("E21").Value = "  +2.01%  "

# Row 22 - Polygon
$ws.Range("E22").Value = "  +2.39%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  +3.49%  "

# Row 24 & 25 swap: Litecoin <-> RenderToken
$ws.Range("B24").Value = "RenderToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D24" "11.55"
$ws.Range("E24").Value = "  +10.65%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell "D25" "83.33"
$ws.Range("E25").Value = "  +1.92%  "

# Row 26 - Fetch.AI
$ws.Range("E26").Value = "  +8.47%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextCell "D27" "12.44"
$ws.Range("E27").Value = "  +4.28%  "

# Row 28 - Dai
$ws.Range("E28").Value = "  -0.04%  "

# Row 29 - NEARProtocol
Set-TextCell "D29" "7.63"
$ws.Range("E29").Value = "  +5.91%  "

# Row 30 - ImmutableX
Set-TextCell "D30" "2.28"
$ws.Range("E30").Value = "  +11.49%  "

# Row 31 - PancakeSwap
Set-TextCell "D31" "2.72"
$ws.Range("E31").Value = "  +1.44%  "

# Row 32 - FirstDigitalUSD
Set-TextCell "D32" "1.00"
$ws.Range("E32").Value = "  +0.07%  "

# Row 33 - EthereumClassic
Set-TextCell "D33" "27.74"
$ws.Range("E33").Value = "  +2.07%  "

# Row 34 - Hedera
Set-TextCell "D34" "0.111"
$ws.Range("E34").Value = "  +2.57%  "

# Row 35 - PEPE
Set-TextCell "D35" "0.0₃0871"
$ws.Range("E35").Value = "  +7.34%  "

# Row 36 - Mantle
$ws.Range("E36").Value = "  +3.19%  "

# Row 37 - Filecoin
Set-TextCell "D37" "5.91"
$ws.Range("E37").Value = "  +2.47%  "

# Row 38 - dogwifhat
Set-TextCell "D38" "3.17"
$ws.Range("E38").Value = "  +12.34%  "

# Row 39 & 40 swap: Kaspa <-> Stacks
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell "D39" "2.13"
$ws.Range("E39").Value = "  +3.72%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell "D40" "0.131"
$ws.Range("E40").Value = "  +9.10%  "

# Row 41 - OKB
Set-TextCell "D41" "50.65"

# Row 42 - Cosmos
Set-TextCell "D42" "9.11"
$ws.Range("E42").Value = "  +1.97%  "

# Row 43 - TheGraph
Set-TextCell "D43" "0.313"
$ws.Range("E43").Value = "  +16.39%  "

# Row 44 - Arweave
Set-TextCell "D44" "43.97"
$ws.Range("E44").Value = "  +14.00%  "

# Row 45 - Bittensor
Set-TextCell "D45" "396.90"
$ws.Range("E45").Value = "  +3.25%  "

# Row 46 - VeChain
$ws.Range("E46").Value = "  +3.64%  "

# Row 47 - Maker
Set-TextCell "D47" "2.729.34"
$ws.Range("E47").Value = "  +1.21%  "

# Row 48 - Monero
Set-TextCell "D48" "133.45"
$ws.Range("E48").Value = "  +2.35%  "

# Row 49 - InjectiveProtocol
Set-TextCell "D49" "26.05"
$ws.Range("E49").Value = "  +12.29%  "

# Row 51 - ThetaToken
$ws.Range("E51").Value = "  +7.22%  "
